$d = $word.ActiveDocument
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "Joel*") {
        $r = $p.Range
        $r.End = $r.End - 1
        $r.Collapse(0)
        $r.InsertAfter("Escobar - 44107580")
        $newRange = $p.Range
        $newRange.End = $newRange.End - 1
        $newRange.Start = $newRange.End - ("Escobar - 44107580".Length)
        # force distinct formatting first
        $newRange.Font.Size = 1
        $newRange.Font.Name = "Arial"
        $newRange.Font.NameBi = "Arial"
        $newRange.Font.Color = 0
        # Now set to match target
        $newRange.Font.Size = 14
        $newRange.Font.SizeBi = 14
        Write-Output "Para text after: [$($p.Range.Text)]"
    }
}
